# Add "contract case" columns to the companies sheet:
# two new boolean-ish indicator columns (isLegalPerson / isNaturalPerson)
# inserted between companyName and minVolume, each row marked with "x"
# for whichever case applies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before the existing "minVolume" column (C),
# which pushes minVolume from C to E and carries column B's cell style
# (s="1") onto the new C/D cells in rows 2 and 4.
$ws.Range("C1:D1").EntireColumn.Insert()

# New header row: isLegalPerson / isNaturalPerson
$ws.Range("C1").Value = "isLegalPerson"
$ws.Range("D1").Value = "isNaturalPerson"

# Mark which case applies per company:
#  - row 2 (JONESBORO ...)  -> legal person (left blank / unmarked)
#  - row 3 (KODREWEX ...)   -> natural person
#  - row 4 (PALOMINO ...)   -> legal person
$ws.Range("D3").Value = "x"
$ws.Range("C4").Value = "x"

# Give the two new columns the same width as column B.
$ws.Range("C1:D1").EntireColumn.ColumnWidth = $ws.Range("B1").EntireColumn.ColumnWidth
